$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force numeric-looking strings to stay as text
# (NumberFormat "@" on a helper cell, then Copy + PasteSpecial values-only
# onto the target so the target cell keeps its original General format).
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "29.949.36"
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  +0.48%  "

$scratch.Value = "1.892.83"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  +0.01%  "

$scratch.Value = "1.000"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.04%  "

$scratch.Value = "0.7785"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.28%  "

$scratch.Value = "243.93"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.04%  "

$scratch.Value = "1.000"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +0.02%  "

$scratch.Value = "0.3133"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.40%  "

$scratch.Value = "25.83"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +2.28%  "

$scratch.Value = "0.07264"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.36%  "

$scratch.Value = "0.08706"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +7.85%  "

$scratch.Value = "2.038.20"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +6.37%  "

$scratch.Value = "0.7738"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +1.09%  "

$scratch.Value = "5.407"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.85%  "

$scratch.Value = "94.46"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +2.40%  "

$scratch.Value = "6.203"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +0.60%  "

$scratch.Value = "29.945.99"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +0.53%  "

$scratch.Value = "13.92"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.22%  "

$scratch.Value = "2.320.11"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +9.89%  "

$scratch.Value = "245.98"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +1.04%  "

$scratch.Value = "0.000007870"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +1.42%  "

$scratch.Value = "8.179"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.95%  "

$scratch.Value = "1.001"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.11%  "

$scratch.Value = "1.001"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.09%  "

$scratch.Value = "0.1671"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +4.25%  "

$scratch.Value = "9.505"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +1.24%  "

$scratch.Value = "162.98"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.70%  "

$scratch.Value = "18.85"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.61%  "

$scratch.Value = "2.053"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.16%  "

$scratch.Value = "1.435"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.00%  "

$scratch.Value = "1.543"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -0.31%  "

$scratch.Value = "4.517"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +1.16%  "

$scratch.Value = "4.125"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +0.65%  "

$scratch.Value = "0.05485"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -0.75%  "

$scratch.Value = "1.247"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -1.19%  "

$scratch.Value = "0.7555"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +1.44%  "

$scratch.Value = "1.004"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +0.67%  "

$scratch.Value = "2.686"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +2.51%  "

$scratch.Value = "0.01960"
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +2.50%  "

$scratch.Value = "2.791"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +0.52%  "

$scratch.Value = "0.4508"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +2.03%  "

$scratch.Value = "74.14"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.51%  "

$scratch.Value = "1.109.14"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -2.44%  "

$scratch.Value = "6.087"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +4.08%  "

$scratch.Value = "0.8528"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.09%  "

$scratch.Value = "2.201.18"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +9.60%  "

$scratch.Value = "1.000"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.01%  "

$scratch.Value = "103.31"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -0.46%  "

$scratch.Value = "1.882"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -0.27%  "

$scratch.Value = "7.596"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +2.12%  "

$scratch.Value = "9.862"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.87%  "

# Rows 42/43: Aave and Maker swapped position in the ranking
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

# Clean up the scratch cell so it leaves no trace in the sheet
$scratch.Clear()
